$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 302, shifting ADL..PPT down by one (ADL was at 302, now at 303)
$ws.Rows("302:302").Insert()

# Copy the formatting of the (now shifted) ADL row's A cell onto the new TPA row's A cell
$ws.Range("A303").Copy()
$ws.Range("A302").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A302").Value = "TPA"
$ws.Range("B302").Value = "Tampa, United States"
$ws.Range("C302").Value = "TPA"
$ws.Range("D302").Value = 27.9755001068
$ws.Range("E302").Value = -82.533203125
$ws.Range("F302").Value = "US"
$ws.Range("G302").Value = "North America"
$ws.Range("H302").Value = "Tampa"
